$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inst1")

# Update the "ini_2" row (row 7) with the new instance values
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 11
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 12
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 11
$ws.Range("K7").Value = 12

# Move the active selection to match the saved view state
$ws.Range("E12").Select()
